$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7730.574239492416
$ws.Range("B2").Value = 128.8429039915403
$ws.Range("C2").Value = 772.5703779935836
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 0.07761389793462425
